$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay as TEXT (matching style s=8,
# same as the rest of the sheet). A plain .Value assignment would be coerced
# to a number, so: set the value first, then re-apply the original cell
# format (copy format only from a same-styled neighbor) to keep it textual
# without altering the number format / style id.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 22.04.2024"

# Row 6
$ws.Range("B6").Value = "26.04."
$ws.Range("C6").Value = "27.04."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU WMKYRX"
$ws.Range("E6").Value = "140,19-"

# Row 7
$ws.Range("B7").Value = "28.04."
$ws.Range("C7").Value = "29.04."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "42,30-"

# Row 8
$ws.Range("B8").Value = "30.04."
$ws.Range("C8").Value = "01.05."
$ws.Range("D8").Value = "KARTENZ./30.04 LIDL RO"
$ws.Range("E8").Value = "45,26-"

# Row 9
$ws.Range("B9").Value = "02.05."
$ws.Range("C9").Value = "03.05."
$ws.Range("D9").Value = "KARTENZ./02.05 EDEKA RO"
$ws.Range("E9").Value = "23,18-"

# Row 10 (previously empty, now populated). E10 changes style from s=12 to
# s=17 (the right-aligned amount style shared with E6:E9,E12). E10 is the
# left cell of the merged range E10:F10, so copying/pasting the *full*
# two-column range (or re-Merge()-ing) would materialize a spurious F10
# cell that isn't in the target. Paste *formats only* onto the single E10
# cell instead: this forks a proper per-cell style (doesn't mutate the
# shared style 12 that E11 still uses) and leaves the merge + F10 alone.
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = "06.05."
$ws.Range("C10").Value = "07.05."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 88517325"
$ws.Range("E10").Value = "38,16-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 11.05.2024"
$ws.Range("E12").Value = "289,09-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 19.05.2024"
